$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data rows (header + 12 data rows)
$data = @(
    @("COMMENTS", "OFFENSIVE (Y or N)"),
    @("bobo tarantado 💩", "Y"),
    @("gago punyeta 💩", "Y"),
    @("puta fuck gago 💩", "Y"),
    @("maganda", "N"),
    @("tae bobo", "Y"),
    @("pogi", "N"),
    @("matalino mabait", "N"),
    @("masipag magalang atin", "N"),
    @("Yan ung sunod na magdadala Ng pandemia sa [NAME]😠😠😠", "N"),
    @("Naawa ako sa bata at sa magulang niya🥺🥺🥺🥺🥺🥺", "N"),
    @("Nakakatakot naman ang ginawa ni tatay at dumugo pa ang kamay nya fuck🖕", "Y"),
    @("[NAME] bless kuya [NAME] ♥️♥️", $null)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    if ($data[$i][1] -ne $null) {
        $ws.Cells.Item($r, 2).Value = $data[$i][1]
    }
}

# Apply the same center-alignment style as the existing B column cells
# to the newly added rows' B cells (including the empty B13).
$ws.Range("B10:B13").HorizontalAlignment = -4108  # xlCenter

# Column widths (column A is widened to fit the new long comment text; column B
# keeps its original width untouched). The host snaps ColumnWidth to pixel
# boundaries like real Excel, so feed it the pre-image that lands closest to
# the target stored width of 138.42578125.
$ws.Columns.Item(1).ColumnWidth = 137.6666666666667

# Selection
$ws.Range("A4").Select()
